$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$vtab = [char]11

$cell = $tbl.Cell(1, 1)
$cell.Range.Text = "38 x 95" + $vtab + "  9    5" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "8|    |"

$cell = $tbl.Cell(1, 2)
$cell.Range.Text = "72 x 28" + $vtab + "  2    8" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "2|    |"

$cell = $tbl.Cell(1, 3)
$cell.Range.Text = "56 x 16" + $vtab + "  1    6" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "6|    |"

$cell = $tbl.Cell(2, 1)
$cell.Range.Text = "96 x 61" + $vtab + "  6    1" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "6|    |"

$cell = $tbl.Cell(2, 2)
$cell.Range.Text = "54 x 93" + $vtab + "  9    3" + $vtab + "  ----" + $vtab + "5|    |" + $vtab + "4|    |"

$cell = $tbl.Cell(2, 3)
$cell.Range.Text = "70 x 11" + $vtab + "  1    1" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "0|    |"

$cell = $tbl.Cell(3, 1)
$cell.Range.Text = "79 x 30" + $vtab + "  3    0" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "9|    |"

$cell = $tbl.Cell(3, 2)
$cell.Range.Text = "62 x 18" + $vtab + "  1    8" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "2|    |"

$cell = $tbl.Cell(3, 3)
$cell.Range.Text = "63 x 99" + $vtab + "  9    9" + $vtab + "  ----" + $vtab + "6|    |" + $vtab + "3|    |"

$cell = $tbl.Cell(4, 1)
$cell.Range.Text = "43 x 60" + $vtab + "  6    0" + $vtab + "  ----" + $vtab + "4|    |" + $vtab + "3|    |"

$cell = $tbl.Cell(4, 2)
$cell.Range.Text = "95 x 88" + $vtab + "  8    8" + $vtab + "  ----" + $vtab + "9|    |" + $vtab + "5|    |"

$cell = $tbl.Cell(4, 3)
$cell.Range.Text = "31 x 15" + $vtab + "  1    5" + $vtab + "  ----" + $vtab + "3|    |" + $vtab + "1|    |"

$cell = $tbl.Cell(5, 1)
$cell.Range.Text = "70 x 49" + $vtab + "  4    9" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "0|    |"

$cell = $tbl.Cell(5, 2)
$cell.Range.Text = "82 x 32" + $vtab + "  3    2" + $vtab + "  ----" + $vtab + "8|    |" + $vtab + "2|    |"

$cell = $tbl.Cell(5, 3)
$cell.Range.Text = "79 x 69" + $vtab + "  6    9" + $vtab + "  ----" + $vtab + "7|    |" + $vtab + "9|    |"
